# algorithm_list.xlsx - "removed keyword scoring, updated keyword generation,
# updated name generation"
#
# Sheet "algorithms": rows 14 downward hold generated keyword-combination
# records (columns B/D/F = keyword_type_1/2/3, C/E/G = modifier_1/2/3,
# H = deactivate flag). This edit:
#  - turns row 14/15 into 2-field "pref"/"suff" records (drops their
#    keyword_type_3 / modifier_3 / deactivate columns, F/G/H)
#  - re-generates the noun/verb/adje permutation table for rows 16-40
#  - appends two more permutation rows (41, 42) so the noun/verb/adje x
#    noun/verb/adje x noun/verb/adje grid (27 rows) is complete

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("algorithms")

# --- row 14: becomes a 2-column "pref" keyword record -----------------
$ws.Range("B14").Value = "pref"
$ws.Range("F14").Clear()
$ws.Range("G14").Clear()
$ws.Range("H14").Clear()

# --- row 15: becomes a 2-column "suff" keyword record ------------------
$ws.Range("D15").Value = "suff"
$ws.Range("F15").Clear()
$ws.Range("G15").Clear()
$ws.Range("H15").Clear()

# --- rows 16-40: regenerated noun/verb/adje permutations ---------------
$ws.Range("D16").Value = "noun"

$ws.Range("D17").Value = "verb"
$ws.Range("F17").Value = "noun"

$ws.Range("D18").Value = "adje"
$ws.Range("F18").Value = "noun"

$ws.Range("D19").Value = "noun"

$ws.Range("D20").Value = "verb"
$ws.Range("F20").Value = "verb"

$ws.Range("D21").Value = "adje"
$ws.Range("F21").Value = "verb"

$ws.Range("D22").Value = "noun"

$ws.Range("B23").Value = "noun"
$ws.Range("D23").Value = "verb"
$ws.Range("F23").Value = "adje"

$ws.Range("B24").Value = "noun"
$ws.Range("D24").Value = "adje"
$ws.Range("F24").Value = "adje"

$ws.Range("D25").Value = "noun"

$ws.Range("D26").Value = "verb"
$ws.Range("F26").Value = "noun"

$ws.Range("D27").Value = "adje"
$ws.Range("F27").Value = "noun"

$ws.Range("D28").Value = "noun"

$ws.Range("D29").Value = "verb"
$ws.Range("F29").Value = "verb"

$ws.Range("D30").Value = "adje"
$ws.Range("F30").Value = "verb"

$ws.Range("D31").Value = "noun"

$ws.Range("B32").Value = "verb"
$ws.Range("D32").Value = "verb"
$ws.Range("F32").Value = "adje"

$ws.Range("B33").Value = "verb"
$ws.Range("D33").Value = "adje"
$ws.Range("F33").Value = "adje"

$ws.Range("D34").Value = "noun"

$ws.Range("D35").Value = "verb"
$ws.Range("F35").Value = "noun"

$ws.Range("D36").Value = "adje"
$ws.Range("F36").Value = "noun"

$ws.Range("D37").Value = "noun"

$ws.Range("D38").Value = "verb"
$ws.Range("F38").Value = "verb"

$ws.Range("D39").Value = "adje"
$ws.Range("F39").Value = "verb"

$ws.Range("D40").Value = "noun"

# --- rows 41-42: two additional permutation rows (new) -----------------
$ws.Range("A41").Value = 39
$ws.Range("B41").Value = "adje"
$ws.Range("C41").Value = "no_cut"
$ws.Range("D41").Value = "verb"
$ws.Range("E41").Value = "no_cut"
$ws.Range("F41").Value = "adje"
$ws.Range("G41").Value = "no_cut"
$ws.Range("H41").Value = "yes"

$ws.Range("A42").Value = 40
$ws.Range("B42").Value = "adje"
$ws.Range("C42").Value = "no_cut"
$ws.Range("D42").Value = "adje"
$ws.Range("E42").Value = "no_cut"
$ws.Range("F42").Value = "adje"
$ws.Range("G42").Value = "no_cut"
$ws.Range("H42").Value = "yes"

# --- view state: match the selection left by the edit ------------------
$ws.Range("A44").Select()
